$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this market/product.
# It belongs right after the existing row 5 (chronologically at the top of
# the data block), so insert a fresh row at position 6 and push every
# following row down by one (old row 6 -> 7, ..., old row 61 -> 62).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly observation.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44552
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112030
$ws.Range("G6").Value = "Poroto granado"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 37000
$ws.Range("L6").Value = 38000
$ws.Range("M6").Value = 37500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1500
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
